$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices stored as literal text (e.g. "62.655.24"); force
# text format before writing so numeric-looking values ("571.62", "0.530",
# etc.) are not auto-converted to floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "62.655.24"
$ws.Range("E2").Value = "  -0.96%  "

$ws.Range("D3").Value = "2.455.06"
$ws.Range("E3").Value = "  -1.11%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "571.62"
$ws.Range("E5").Value = "  -1.18%  "

$ws.Range("D6").Value = "147.62"
$ws.Range("E6").Value = "  +0.42%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "0.530"
$ws.Range("E8").Value = "  -1.64%  "

$ws.Range("D9").Value = "0.110"
$ws.Range("E9").Value = "  -1.46%  "

$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  -0.32%  "

$ws.Range("D11").Value = "5.17"
$ws.Range("E11").Value = "  -1.55%  "

$ws.Range("D12").Value = "0.346"
$ws.Range("E12").Value = "  -2.16%  "

$ws.Range("D13").Value = "28.80"
$ws.Range("E13").Value = "  +0.67%  "

$ws.Range("D14").Value = "0.0000175"
$ws.Range("E14").Value = "  -2.44%  "

$ws.Range("D15").Value = "2.902.88"
$ws.Range("E15").Value = "  -0.86%  "

$ws.Range("D16").Value = "62.562.13"
$ws.Range("E16").Value = "  -0.81%  "

$ws.Range("D17").Value = "2.460.73"
$ws.Range("E17").Value = "  -0.51%  "

$ws.Range("D18").Value = "7.65"
$ws.Range("E18").Value = "  -5.46%  "

$ws.Range("D19").Value = "10.67"
$ws.Range("E19").Value = "  -3.30%  "

$ws.Range("D20").Value = "2.31"
$ws.Range("E20").Value = "  +1.35%  "

$ws.Range("E21").Value = "  +0.49%  "

$ws.Range("D22").Value = "320.94"
$ws.Range("E22").Value = "  -2.60%  "

$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").Value = "10.34"
$ws.Range("E24").Value = "  +4.34%  "

$ws.Range("D25").Value = "64.55"
$ws.Range("E25").Value = "  -2.52%  "

$ws.Range("D26").Value = "637.33"
$ws.Range("E26").Value = "  -3.51%  "

$ws.Range("D27").Value = "2.572.42"
$ws.Range("E27").Value = "  -1.03%  "

$ws.Range("D28").Value = "0.0₃0963"
$ws.Range("E28").Value = "  -3.41%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("D30").Value = "1.41"
$ws.Range("E30").Value = "  -4.68%  "

$ws.Range("D31").Value = "7.86"
$ws.Range("E31").Value = "  -2.98%  "

$ws.Range("D32").Value = "1.81"
$ws.Range("E32").Value = "  -2.86%  "

$ws.Range("D33").Value = "0.132"
$ws.Range("E33").Value = "  -0.89%  "

$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.06%  "

$ws.Range("D35").Value = "1.49"
$ws.Range("E35").Value = "  -4.03%  "

$ws.Range("D36").Value = "4.67"
$ws.Range("E36").Value = "  -2.59%  "

$ws.Range("D37").Value = "5.38"
$ws.Range("E37").Value = "  -1.49%  "

$ws.Range("D38").Value = "0.365"
$ws.Range("E38").Value = "  -1.92%  "

$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "149.84"
$ws.Range("E39").Value = "  -0.81%  "

$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").Value = "18.49"
$ws.Range("E40").Value = "  -1.76%  "

$ws.Range("E41").Value = "  -2.48%  "

$ws.Range("D42").Value = "1.72"
$ws.Range("E42").Value = "  -2.03%  "

$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "1.01"
$ws.Range("E43").Value = "  +0.72%  "

$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").Value = "0.0₆0299"
$ws.Range("E44").Value = "  -7.65%  "

$ws.Range("D45").Value = "154.04"
$ws.Range("E45").Value = "  -0.51%  "

$ws.Range("D46").Value = "15.39"
$ws.Range("E46").Value = "  +0.81%  "

$ws.Range("D47").Value = "3.55"
$ws.Range("E47").Value = "  -1.84%  "

$ws.Range("D48").Value = "0.605"
$ws.Range("E48").Value = "  -0.85%  "

$ws.Range("D49").Value = "20.19"
$ws.Range("E49").Value = "  -1.20%  "

$ws.Range("E50").Value = "  -0.94%  "

$ws.Range("D51").Value = "0.0900"
$ws.Range("E51").Value = "  -2.00%  "

# Restore the original (default) cell formatting now that the literal text
# has been stored, so the saved style matches the unformatted source cells.
$ws.Range("D2:D51").ClearFormats()